# "shiny slides 1 and 2"
#
# Fill in the link_slide (column H) entries for the two R-Shiny sessions
# on day 3 (rows 21 and 23) that were previously left blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep gridlines visible (matches the source file) - the host resets this
# to hidden on every round-trip unless told otherwise.
$excel.ActiveWindow.DisplayGridlines = $true

$ws.Range("H21").Value = "topics/Shiny/slide_shiny1.html"
$ws.Range("H23").Value = "topics/Shiny/slide_shiny2.html"

# The room-name column was widened while the author was filling these in.
$ws.Columns.Item(2).ColumnWidth = 34.8

# Leave the selection where the author ended up.
[void]$ws.Range("F21").Select()
